$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has a "Programa resumido:" entry at row 12.
# We need to insert three new rows above it to hold a new
# "Docentes responsaveis:" (Responsible faculty) section:
#   row 12 -> "Docentes responsaveis:" label (column A only)
#   row 13 -> first professor's name (columns B and C only)
#   row 14 -> second professor's name (columns B and C only)
# Inserting at row 12 three times pushes the existing rows 12-20 down
# to rows 15-23, carrying their values, shared-string type and styles
# with them automatically.
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()

# Row 12: label only in column A (bold style already applied by Insert)
$ws.Range("A12").Value2 = "Docentes responsáveis:"
$ws.Range("B12").Clear()
$ws.Range("C12").Clear()

# Row 13: first responsible professor, columns B and C only
$ws.Range("A13").Clear()
$ws.Range("B13").Value2 = "198273 - Domingos Savio Giordani"
$ws.Range("C13").Value2 = "198273 - Domingos Savio Giordani"

# Row 14: second responsible professor, columns B and C only
$ws.Range("A14").Clear()
$ws.Range("B14").Value2 = "1506103 - Pedro Carlos de Oliveira"
$ws.Range("C14").Value2 = "1506103 - Pedro Carlos de Oliveira"

# These three new rows use the default (automatic) row height, unlike
# the custom heights used by the surrounding rows.
$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()
